$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1.006090927296766
$ws.Range("D2").Value = 1.008736917786339
$ws.Range("E2").Value = 1.008722147100834
$ws.Range("F2").Value = 1.004255285924022
$ws.Range("J2").Value = 1.011370418211365
$ws.Range("K2").Value = 1.011612260904118
$ws.Range("L2").Value = 1.011597535032606
$ws.Range("M2").Value = 1.00714429220932
$ws.Range("N2").Value = 1.007347764093004
$ws.Range("C3").Value = 1.008747479273434
$ws.Range("D3").Value = 1.011303855525857
$ws.Range("E3").Value = 1.011042418778869
$ws.Range("F3").Value = 1.007591720365493
$ws.Range("J3").Value = 1.013648668426333
$ws.Range("K3").Value = 1.013977948692831
$ws.Range("L3").Value = 1.013717244929004
$ws.Range("M3").Value = 1.010276259819716
$ws.Range("N3").Value = 1.008163843152681
$ws.Range("C4").Value = 1.010454205672261
$ws.Range("D4").Value = 1.01295315082762
$ws.Range("E4").Value = 1.012532685184043
$ws.Range("F4").Value = 1.009736437635741
$ws.Range("J4").Value = 1.015110944212857
$ws.Range("K4").Value = 1.015496830735572
$ws.Range("L4").Value = 1.015077481814617
$ws.Range("M4").Value = 1.012288686452855
$ws.Range("N4").Value = 1.008686106013879
$ws.Range("C5").Value = 1.011168864086609
$ws.Range("D5").Value = 1.013643794106838
$ws.Range("E5").Value = 1.013156603907739
$ws.Range("F5").Value = 1.010634790494186
$ws.Range("J5").Value = 1.015722906349844
$ws.Range("K5").Value = 1.016132598703524
$ws.Range("L5").Value = 1.015646672322156
$ws.Range("M5").Value = 1.013131421759703
$ws.Range("N5").Value = 1.008904303239479
$ws.Range("C6").Value = 1.011288693616128
$ws.Range("D6").Value = 1.013759598832089
$ws.Range("E6").Value = 1.013261212756075
$ws.Range("F6").Value = 1.010785438287012
$ws.Range("J6").Value = 1.015825496594371
$ws.Range("K6").Value = 1.016239186605167
$ws.Range("L6").Value = 1.015742088224938
$ws.Range("M6").Value = 1.013272730806825
$ws.Range("N6").Value = 1.008940860364173
$ws.Range("C7").Value = 1.010463766062007
$ws.Range("D7").Value = 1.012962389824795
$ws.Range("E7").Value = 1.012541032099204
$ws.Range("F7").Value = 1.00974845424076
$ws.Range("J7").Value = 1.015119132104623
$ws.Range("K7").Value = 1.015505336686019
$ws.Range("L7").Value = 1.015085097703791
$ws.Range("M7").Value = 1.012299959911009
$ws.Range("N7").Value = 1.008689026894543
$ws.Range("C8").Value = 1.00699130524273
$ws.Range("D8").Value = 1.009606892739671
$ws.Range("E8").Value = 1.009508635157861
$ws.Range("F8").Value = 1.005385851280336
$ws.Range("J8").Value = 1.012142871873044
$ws.Range("K8").Value = 1.01241425972983
$ws.Range("L8").Value = 1.012316292585387
$ws.Range("M8").Value = 1.008205750865688
$ws.Range("N8").Value = 1.007624777058741
$ws.Range("C9").Value = 1.000774926865739
$ws.Range("D9").Value = 1.00360105765439
$ws.Range("E9").Value = 1.004076946789064
$ws.Range("F9").Value = 0.997584833785194
$ws.Range("J9").Value = 1.006803877021899
$ws.Range("K9").Value = 1.006873078252644
$ws.Range("L9").Value = 1.007347303589048
$ws.Range("M9").Value = 1.000878039125616
$ws.Range("N9").Value = 1.005703893909823
$ws.Range("C10").Value = 0.9965598517892567
$ws.Range("D10").Value = 0.9995295930248154
$ws.Range("E10").Value = 1.000391984513437
$ws.Range("F10").Value = 0.9923007265657137
$ws.Range("J10").Value = 1.00317638050384
$ws.Range("K10").Value = 1.003110765348437
$ws.Range("L10").Value = 1.003969824468227
$ws.Range("M10").Value = 0.9959100850204848
$ws.Range("N10").Value = 1.004391033268429
$ws.Range("C11").Value = 0.9947166684116088
$ws.Range("D11").Value = 0.9977494253751157
$ws.Range("E11").Value = 0.9987801794014731
$ws.Range("F11").Value = 0.9899912410929472
$ws.Range("J11").Value = 1.001588398655513
$ws.Range("K11").Value = 1.001464374056901
$ws.Range("L11").Value = 1.002490978865181
$ws.Range("M11").Value = 0.9937377287925286
$ws.Range("N11").Value = 1.003814516701121
$ws.Range("C12").Value = 0.9940292131350527
$ws.Range("D12").Value = 0.9970855071619035
$ws.Range("E12").Value = 0.9981789596307084
$ws.Range("F12").Value = 0.9891300330018744
$ws.Range("J12").Value = 1.000995865763107
$ws.Range("K12").Value = 1.000850138015531
$ws.Range("L12").Value = 1.00193912301098
$ws.Range("M12").Value = 0.9929274990880825
$ws.Range("N12").Value = 1.003599131600842
$ws.Range("C13").Value = 0.9941768037503463
$ws.Range("D13").Value = 0.9972280430177404
$ws.Range("E13").Value = 0.9983080390005992
$ws.Range("F13").Value = 0.9893149195511527
$ws.Range("J13").Value = 1.001123089153445
$ws.Range("K13").Value = 1.000982017172393
$ws.Range("L13").Value = 1.002057614631271
$ws.Range("M13").Value = 0.9931014485550306
$ws.Range("N13").Value = 1.00364538914935
$ws.Range("C14").Value = 0.9946599011945081
$ws.Range("D14").Value = 0.9976946010513851
$ws.Range("E14").Value = 0.9987305344040769
$ws.Range("F14").Value = 0.989920122794787
$ws.Range("J14").Value = 1.001539475043138
$ws.Range("K14").Value = 1.001413656617626
$ws.Range("L14").Value = 1.00244541475185
$ws.Range("M14").Value = 0.9936708234728931
$ws.Range("N14").Value = 1.003796738418048
$ws.Range("C15").Value = 0.9949571772351183
$ws.Range("D15").Value = 0.9979817040314933
$ws.Range("E15").Value = 0.998990510600829
$ws.Range("F15").Value = 0.990292558607223
$ws.Range("J15").Value = 1.001795665392203
$ws.Range("K15").Value = 1.001679244150521
$ws.Range("L15").Value = 1.002684011100384
$ws.Range("M15").Value = 0.9940211901285849
$ws.Range("N15").Value = 1.003889824213849
$ws.Range("C16").Value = 0.9966817889603924
$ws.Range("D16").Value = 0.999647366016141
$ws.Range("E16").Value = 1.000498605863378
$ws.Range("F16").Value = 0.9924535359026349
$ws.Range("J16").Value = 1.003281398206677
$ws.Range("K16").Value = 1.003219658601092
$ws.Range("L16").Value = 1.004067618246911
$ws.Range("M16").Value = 0.9960537990814301
$ws.Range("N16").Value = 1.004429122571202
$ws.Range("C17").Value = 0.9977586915948781
$ws.Range("D17").Value = 1.000687517331306
$ws.Range("E17").Value = 1.001440195758832
$ws.Range("F17").Value = 0.9938032226934108
$ws.Range("J17").Value = 1.004208674963362
$ws.Range("K17").Value = 1.004181225239359
$ws.Range("L17").Value = 1.004931074105202
$ws.Range("M17").Value = 0.9973230313618571
$ws.Range("N17").Value = 1.004765234606173
$ws.Range("C18").Value = 0.9983850985791192
$ws.Range("D18").Value = 1.0012925677505
$ws.Range("E18").Value = 1.001987853182728
$ws.Range("F18").Value = 0.9945884141047743
$ws.Range("J18").Value = 1.00474788172371
$ws.Range("K18").Value = 1.004740429677373
$ws.Range("L18").Value = 1.005433139341394
$ws.Range("M18").Value = 0.9980613169681083
$ws.Range("N18").Value = 1.004960509703678
$ws.Range("C19").Value = 0.9985983965555357
$ws.Range("D19").Value = 1.001498597077151
$ws.Range("E19").Value = 1.002174328996302
$ws.Range("F19").Value = 0.9948557992536934
$ws.Range("J19").Value = 1.00493145885666
$ws.Range("K19").Value = 1.004930825056835
$ws.Range("L19").Value = 1.005604066228444
$ws.Range("M19").Value = 0.9983127119490669
$ws.Range("N19").Value = 1.005026963285529
$ws.Range("C20").Value = 0.9976433300527123
$ws.Range("D20").Value = 1.000576090560031
$ws.Range("E20").Value = 1.001339333670625
$ws.Range("F20").Value = 0.993658627964647
$ws.Range("J20").Value = 1.004109359115974
$ws.Range("K20").Value = 1.004078230742143
$ws.Range("L20").Value = 1.004838596883237
$ws.Range("M20").Value = 0.9971870663168665
$ws.Range("N20").Value = 1.004729253196774
$ws.Range("C21").Value = 0.9945177196130314
$ws.Range("D21").Value = 0.9975572862761313
$ws.Range("E21").Value = 0.9986061904447127
$ws.Range("F21").Value = 0.9897419995063174
$ws.Range("J21").Value = 1.001416934683185
$ws.Range("K21").Value = 1.001286624685652
$ws.Range("L21").Value = 1.002331288293003
$ws.Range("M21").Value = 0.9935032494993649
$ws.Range("N21").Value = 1.003752204366218
$ws.Range("C22").Value = 0.9925361831080421
$ws.Range("D22").Value = 0.9956436591728185
$ws.Range("E22").Value = 0.9968731071297999
$ws.Range("F22").Value = 0.9872599324566337
$ws.Range("J22").Value = 0.9997085132618067
$ws.Range("K22").Value = 0.9995158005207417
$ws.Range("L22").Value = 1.000740064360187
$ws.Range("M22").Value = 0.9911678100289042
$ws.Range("N22").Value = 1.003130696647749
$ws.Range("C23").Value = 0.9935882184135094
$ws.Range("D23").Value = 0.9966596209353727
$ws.Range("E23").Value = 0.9977932668817237
$ws.Range("F23").Value = 0.9885776228259336
$ws.Range("J23").Value = 1.000615689345495
$ws.Range("K23").Value = 1.000456062375759
$ws.Range("L23").Value = 1.001585032801825
$ws.Range("M23").Value = 0.9924077441135356
$ws.Range("N23").Value = 1.003460863377317
$ws.Range("C24").Value = 0.9976954623140011
$ws.Range("D24").Value = 1.000626444621296
$ws.Range("E24").Value = 1.001384913709678
$ws.Range("F24").Value = 0.993723970449555
$ws.Range("J24").Value = 1.00415424079048
$ws.Range("K24").Value = 1.004124774649397
$ws.Range("L24").Value = 1.004880388217769
$ws.Range("M24").Value = 0.9972485093517661
$ws.Range("N24").Value = 1.004745514034718
$ws.Range("C25").Value = 1.002394095298738
$ws.Range("D25").Value = 1.005165242897434
$ws.Range("E25").Value = 1.00549208066899
$ws.Range("F25").Value = 0.9996157571292116
$ws.Range("J25").Value = 1.008195799769047
$ws.Range("K25").Value = 1.008317266850159
$ws.Range("L25").Value = 1.008643008013573
$ws.Range("M25").Value = 1.002786521031373
$ws.Range("N25").Value = 1.006206052205763
